# Applies the price/volume refresh + the FraxShare/WEMIXToken row swap
# described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "26.950.70"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "'" + "1.556.45"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'" + "206.75"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'" + "0.487"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'" + "0.0858"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'" + "1.778.52"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'" + "1.557.32"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'" + "26.938.30"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "'" + "61.71"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "'" + "214.70"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'" + "0.0₃0688"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "'" + "9.20"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").Value = "'" + "153.42"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").Value = "'" + "14.89"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").Value = "'" + "1.09"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").Value = "'" + "1.370.07"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "'" + "0.973"
$ws.Range("E36").Value = "  +5.56%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "'" + "0.523"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'" + "0.991"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'" + "5.52"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("D45").Value = "'" + "63.75"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "'" + "1.691.50"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'" + "86.05"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'" + "0.0956"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "'" + "1.00"
$ws.Range("E51").Value = "  +0.32%  "
